$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add the 4 new book rows (12-15) - typed left to right, top to bottom,
#    so new shared strings are appended in that natural order.
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "El que se duerme pierde"
$ws.Cells.Item(12, 3).Value = "Tom Peter"
$ws.Cells.Item(12, 4).Value = 16

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Sin lugar a duda"
$ws.Cells.Item(13, 3).Value = "Ana Gutierrez"
$ws.Cells.Item(13, 4).Value = 26

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "El arte de dormir"
$ws.Cells.Item(14, 3).Value = "Nico"
$ws.Cells.Item(14, 4).Value = 32

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Buscando a Nemo"
$ws.Cells.Item(15, 3).Value = "Humble Po"
$ws.Cells.Item(15, 4).Value = 41

# 2) Update the Author of row 2 to "Max" (new shared string, added last).
$ws.Cells.Item(2, 3).Value = "Max"

# 3) Update the Price of row 6.
$ws.Cells.Item(6, 4).Value = 37090

# 4) Collapse the bestFit columns B and C.
$ws.Columns.Item(2).EntireColumn.Collapsed = $true
$ws.Columns.Item(3).EntireColumn.Collapsed = $true
